# Fruta / hortaliza, semanal
# Weekly refresh: insert the newest price record as a new row 187
# (pushing the existing historical rows 187-278 down to 188-279).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("187:187").Insert()

$ws.Range("A187").Value = 7
$ws.Range("B187").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C187").Value = "Ñuble"
$ws.Range("D187").Value = 44960
$ws.Range("E187").Value = 16
$ws.Range("F187").Value = "Fruta"
$ws.Range("G187").Value = 100108
$ws.Range("H187").Value = "Tropicales y subtropicales"
$ws.Range("I187").Value = 100108005
$ws.Range("J187").Value = "Piña"
$ws.Range("K187").Value = "Caramelo"
$ws.Range("L187").Value = "Segunda"
$ws.Range("M187").Value = 30
$ws.Range("N187").Value = 18000
$ws.Range("O187").Value = 18000
$ws.Range("P187").Value = 18000
$ws.Range("Q187").Value = "$/caja 14 unidades"
$ws.Range("R187").Value = "Ecuador"
$ws.Range("S187").Value = 1286
$ws.Range("T187").Value = 14
